$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ANG001) - D2 already has "Mecanicien"
$ws.Range("C2").Value = "pascadre"

# Row 3 (ANG002)
$ws.Range("C3").Value = "pas cadre"
$ws.Range("D3").Value = "agent de securite"

# Row 4 (ANG003)
$ws.Range("C4").Value = "cadre"
$ws.Range("D4").Value = "Informaticien"

# Row 5 (ANG004)
$ws.Range("C5").Value = " cadre"
$ws.Range("D5").Value = "developpeur"

# Row 6 (ANG005)
$ws.Range("C6").Value = "cadre"
$ws.Range("D6").Value = "developpeur"

# Row 7 (ANG006)
$ws.Range("C7").Value = "pas cadre"
$ws.Range("D7").Value = "Mecanicien"

# Row 8 (ANG007)
$ws.Range("C8").Value = "pas cadre"
$ws.Range("D8").Value = "Mecanicien"

# Row 9 (ANG008)
$ws.Range("C9").Value = "pas cadre"
$ws.Range("D9").Value = "Mecanicien"

# Row 10 (ANG009)
$ws.Range("C10").Value = "pas cadre"
$ws.Range("D10").Value = "Mecanicien"

# Row 11 (ANG010)
$ws.Range("C11").Value = "pas cadre"
$ws.Range("D11").Value = "Mecanicien"

# Row 12 (ANG011)
$ws.Range("C12").Value = "pas cadre"
$ws.Range("D12").Value = "Mecanicien"

# Row 13 (ANG011)
$ws.Range("C13").Value = "pas cadre"
$ws.Range("D13").Value = "Mecanicien"

# Row 14 (ANG012)
$ws.Range("C14").Value = "pas cadre"
$ws.Range("D14").Value = "Mecanicien"

# Row 15 (ANG013)
$ws.Range("C15").Value = "pas cadre"
$ws.Range("D15").Value = "Mecanicien"

# Row 16 (ANG014)
$ws.Range("C16").Value = "pas cadre"
$ws.Range("D16").Value = "Mecanicien"

# Row 17 (ANG015)
$ws.Range("C17").Value = "pas cadre"
$ws.Range("D17").Value = "Mecanicien"

# Row 18 (ANG016)
$ws.Range("C18").Value = "pas cadre"
$ws.Range("D18").Value = "Mecanicien"

# Row 19 (ANG017)
$ws.Range("C19").Value = "cadre"
$ws.Range("D19").Value = "gestionnaire comptable"

# Row 20 (ANG018)
$ws.Range("C20").Value = "pas cadre"
$ws.Range("D20").Value = "agent d entretien"

# Row 21 (ANG019)
$ws.Range("C21").Value = "cadre"
$ws.Range("D21").Value = "gestionnaire production"

# Update column widths (bestFit) to match new content widths
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(4).ColumnWidth = 21.666666666666664

# Update selection to D22 (matches the cell after last edited column)
$ws.Range("D22").Select() | Out-Null
